# Edit: add a new week of price data for Vega Monumental Concepción - Zanahoria.
# This inserts two new rows right before the current row 284, shifting the
# existing rows 284:354 down to 286:356, and fills the two new rows with the
# new week's data (fecha 44995 = 2023-03-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 284 (old row 284 becomes row 286, etc.)
$ws.Range("A284:A285").EntireRow.Insert()

# Fill in the new row 284 (Primera)
$ws.Cells.Item(284, 1).Value  = 11
$ws.Cells.Item(284, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(284, 3).Value  = "Bíobío"
$ws.Cells.Item(284, 4).Value  = 44995
$ws.Cells.Item(284, 5).Value  = 8
$ws.Cells.Item(284, 6).Value  = 100114013
$ws.Cells.Item(284, 7).Value  = "Zanahoria"
$ws.Cells.Item(284, 8).Value  = "Sin especificar"
$ws.Cells.Item(284, 9).Value  = "Primera"
$ws.Cells.Item(284, 10).Value = 800
$ws.Cells.Item(284, 11).Value = 6000
$ws.Cells.Item(284, 12).Value = 6500
$ws.Cells.Item(284, 13).Value = 6250
$ws.Cells.Item(284, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(284, 15).Value = "Región de Ñuble"
$ws.Cells.Item(284, 16).Value = 312
$ws.Cells.Item(284, 17).Value = 20
$ws.Cells.Item(284, 18).Value = "Hortaliza"

# Fill in the new row 285 (Segunda)
$ws.Cells.Item(285, 1).Value  = 11
$ws.Cells.Item(285, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(285, 3).Value  = "Bíobío"
$ws.Cells.Item(285, 4).Value  = 44995
$ws.Cells.Item(285, 5).Value  = 8
$ws.Cells.Item(285, 6).Value  = 100114013
$ws.Cells.Item(285, 7).Value  = "Zanahoria"
$ws.Cells.Item(285, 8).Value  = "Sin especificar"
$ws.Cells.Item(285, 9).Value  = "Segunda"
$ws.Cells.Item(285, 10).Value = 400
$ws.Cells.Item(285, 11).Value = 5000
$ws.Cells.Item(285, 12).Value = 5000
$ws.Cells.Item(285, 13).Value = 5000
$ws.Cells.Item(285, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(285, 15).Value = "Región de Ñuble"
$ws.Cells.Item(285, 16).Value = 250
$ws.Cells.Item(285, 17).Value = 20
$ws.Cells.Item(285, 18).Value = "Hortaliza"
